$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) holds one copy of the rows that changed.
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F17").Value = 1922
$wsExpo.Range("F21").Value = 201
$wsExpo.Range("F24").Value = 174
$wsExpo.Range("F28").Value = 3159
$wsExpo.Range("F30").Value = 115

# Sheet "全部类型" (All types) holds the same events (offset by a few rows)
# and needs the identical updates applied.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F22").Value = 1922
$wsAll.Range("F25").Value = 201
$wsAll.Range("F29").Value = 174
$wsAll.Range("F31").Value = 3159
$wsAll.Range("F33").Value = 115

$wb.Save()
